# Error Calculations and Plots
# This script re-derives the "missing data" selection for the BCDF/20/seed4 sheet:
#  - two rows (RM 232 and SC 92) are removed entirely (shifting later rows up)
#  - a handful of cells switch between having a numeric value and being blank
#    (simulating which values were "removed" for the imputation experiment)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that no longer appear in the data set.
# Delete the lower row first so the upper row index stays valid.
$ws.Rows(28).Delete()
$ws.Rows(26).Delete()

# Cells that now carry a value (previously blank/missing)
$ws.Range("D2").Value = -13.5
$ws.Range("F6").Value = 16.43
$ws.Range("E8").Value = -6.6
$ws.Range("E10").Value = -6.1
$ws.Range("D11").Value = -15.5
$ws.Range("F11").Value = 17.65
$ws.Range("F13").Value = 17.1
$ws.Range("E15").Value = -8.4
$ws.Range("F18").Value = 18.35
$ws.Range("D21").Value = -14.3
$ws.Range("E25").Value = -7.1
$ws.Range("F25").Value = 16.6
$ws.Range("E27").Value = -10
$ws.Range("C33").Value = 10.4
$ws.Range("D33").Value = -14.1

# Cells that are now blank (previously contained a value)
$ws.Range("D3").Value = ""
$ws.Range("D4").Value = ""
$ws.Range("E5").Value = ""
$ws.Range("E12").Value = ""
$ws.Range("F12").Value = ""
$ws.Range("D13").Value = ""
$ws.Range("F17").Value = ""
$ws.Range("E18").Value = ""
$ws.Range("E19").Value = ""
$ws.Range("F19").Value = ""
$ws.Range("F24").Value = ""
$ws.Range("D25").Value = ""
$ws.Range("C29").Value = ""
$ws.Range("E29").Value = ""
$ws.Range("F31").Value = ""
$ws.Range("F32").Value = ""
$ws.Range("E33").Value = ""
